$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the "Nr" column (A30:A39) after a row was effectively removed
# from the earlier numbering sequence (33 -> 29, 34 -> 30, ... 42 -> 38).
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31
$ws.Range("A33").Value = 32
$ws.Range("A34").Value = 33
$ws.Range("A35").Value = 34
$ws.Range("A36").Value = 35
$ws.Range("A37").Value = 36
$ws.Range("A38").Value = 37
$ws.Range("A39").Value = 38

# Extend the numbering column's formatting one row further down (A40),
# matching the formatting already used by the numbering cells above it.
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active cell / selection to reflect where editing continued.
[void]$ws.Range("A40").Select()
